# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the crypto symbol
# list with the latest scrape values (GitHub Actions daily update).
#
# The source sheet stores these columns as plain text (e.g. "321.37",
# "7.98%") rather than numbers, even though the strings look numeric. Excel's
# COM automation auto-converts a numeric-looking string assigned to
# .Value into a real number/percentage, so each target cell is first forced
# to Text format ("@") before the new string is written — this preserves the
# original "text that looks like a number" representation instead of
# silently turning the cell into a numeric type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2";  Value = "321.38" }
    @{ Cell = "E2";  Value = "7.91%" }
    @{ Cell = "D3";  Value = "48.65" }
    @{ Cell = "E3";  Value = "16.28%" }
    @{ Cell = "E4";  Value = "4.63%" }
    @{ Cell = "D5";  Value = "0.08098" }
    @{ Cell = "E5";  Value = "7.64%" }
    @{ Cell = "D6";  Value = "4.599" }
    @{ Cell = "E6";  Value = "5.41%" }
    @{ Cell = "D7";  Value = "1.642" }
    @{ Cell = "E7";  Value = "2.93%" }
    @{ Cell = "D8";  Value = "1.202" }
    @{ Cell = "E8";  Value = "31.37%" }
    @{ Cell = "D9";  Value = "0.1299" }
    @{ Cell = "E9";  Value = "10.41%" }
    @{ Cell = "D10"; Value = "0.1950" }
    @{ Cell = "E10"; Value = "6.72%" }
    @{ Cell = "D11"; Value = "0.09553" }
    @{ Cell = "E11"; Value = "7.95%" }
    @{ Cell = "D12"; Value = "0.04586" }
    @{ Cell = "E12"; Value = "11.13%" }
    @{ Cell = "D13"; Value = "0.1048" }
    @{ Cell = "E13"; Value = "-0.11%" }
    @{ Cell = "D14"; Value = "0.001335" }
    @{ Cell = "E14"; Value = "4.44%" }
    @{ Cell = "D15"; Value = "0.005850" }
    @{ Cell = "E15"; Value = "-2.57%" }
    @{ Cell = "D16"; Value = "3.345" }
    @{ Cell = "E16"; Value = "0.05%" }
    @{ Cell = "D17"; Value = "2.437" }
    @{ Cell = "E17"; Value = "1.51%" }
    @{ Cell = "E18"; Value = "2.13%" }
    @{ Cell = "D19"; Value = "8.217" }
    @{ Cell = "E19"; Value = "-0.88%" }
    @{ Cell = "D20"; Value = "0.1410" }
    @{ Cell = "D21"; Value = "0.3125" }
    @{ Cell = "E21"; Value = "0.65%" }
    @{ Cell = "D22"; Value = "0.04292" }
    @{ Cell = "E22"; Value = "4.82%" }
    @{ Cell = "D23"; Value = "0.001305" }
    @{ Cell = "D24"; Value = "0.004251" }
    @{ Cell = "E24"; Value = "9.33%" }
    @{ Cell = "E25"; Value = "3.68%" }
    @{ Cell = "D26"; Value = "0.0003539" }
    @{ Cell = "E26"; Value = "-4.96%" }
    @{ Cell = "D38"; Value = "0.02666" }
    @{ Cell = "E38"; Value = "11.42%" }
    @{ Cell = "D39"; Value = "0.05642" }
    @{ Cell = "E39"; Value = "8.16%" }
    @{ Cell = "D40"; Value = "0.006300" }
    @{ Cell = "E40"; Value = "-10.38%" }
    @{ Cell = "D41"; Value = "0.007685" }
    @{ Cell = "E41"; Value = "-1.11%" }
    @{ Cell = "D42"; Value = "0.1439" }
    @{ Cell = "E42"; Value = "8.63%" }
    @{ Cell = "E43"; Value = "3.49%" }
    @{ Cell = "E44"; Value = "14.14%" }
    @{ Cell = "E45"; Value = "-1.34%" }
    @{ Cell = "D46"; Value = "0.00006985" }
    @{ Cell = "E46"; Value = "6.02%" }
    @{ Cell = "E47"; Value = "-0.14%" }
    @{ Cell = "D48"; Value = "0.05354" }
    @{ Cell = "E48"; Value = "18.10%" }
    @{ Cell = "D49"; Value = "0.004000" }
    @{ Cell = "E49"; Value = "-4.89%" }
    @{ Cell = "D50"; Value = "0.00002100" }
    @{ Cell = "E50"; Value = "-0.14%" }
    @{ Cell = "D51"; Value = "0.0002000" }
    @{ Cell = "E51"; Value = "-0.14%" }
)

foreach ($update in $updates) {
    $cell = $ws.Range($update.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $update.Value
}

Write-Output "Updated $($updates.Count) cells"
